$wb = $excel.ActiveWorkbook

$wsNP5 = $wb.Worksheets.Item("NP5")
$wsNP7 = $wb.Worksheets.Item("NP7_raw")

# --- Add summary columns (Sample / Geometric Mean Avg. / Stdev) to NP7_raw,
#     mirroring the layout already present on NP5 ---

$labels = @("Sample", "Geometric Mean Avg.", "Stdev")
$wsNP7.Range("F1").Value = $labels[0]
$wsNP7.Range("G1").Value = $labels[1]
$wsNP7.Range("H1").Value = $labels[2]

$groupNames = @("Untreated", "pDNA", "jetPEI", "LPF2000", "S", "B", "G1", "G2", "G3")
$startRows  = @(2, 5, 8, 11, 14, 17, 20, 23, 26)

for ($i = 0; $i -lt $groupNames.Length; $i++) {
    $row = $i + 2
    $start = $startRows[$i]
    $end = $start + 2

    $wsNP7.Range("F$row").Value = $groupNames[$i]
    $wsNP7.Range("G$row").Formula = "=AVERAGE(B$($start):B$($end))"
    $wsNP7.Range("H$row").Formula = "=STDEV(B$($start):B$($end))"
}

# New columns need widths (matching the widths used on the NP5 sheet's
# Sample/Geometric Mean Avg./Stdev columns)
$wsNP7.Columns.Item(7).ColumnWidth = 17.17
$wsNP7.Columns.Item(8).ColumnWidth = 21

# --- Tab / selection swap: NP7_raw becomes the selected/active tab,
#     NP5 is no longer tab-selected (just keeps its own remembered selection) ---

[void]$wsNP5.Range("F1:H10").Select()
[void]$wsNP7.Range("F1:H10").Select()
